$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 19: montly figures 10.000 -> 1.000.000 (B19:L19)
# ---------------------------------------------------------------------
$ws.Range("B19:L19").Value = 1000000

# ---------------------------------------------------------------------
# Row 20: monthly figures 2.000 -> 20.000 (B20:L20)
# ---------------------------------------------------------------------
$ws.Range("B20:L20").Value = 20000

# ---------------------------------------------------------------------
# Row 23 (App Store / iOS block): new hard-coded total + updated count
# ---------------------------------------------------------------------
$ws.Range("A23").Value = 6500000
$ws.Range("B23").Value = 27
# C23 (Cut=3) and D23's formula / E23 ("iOS") are unchanged.

# ---------------------------------------------------------------------
# Row 24 (Google Play block): new hard-coded total + new count/cut cells
# and its own price formula, mirroring row 23.
# ---------------------------------------------------------------------
$ws.Range("A24").Value = 6500000
$ws.Range("B24").Value = 27
$ws.Range("C24").Value = 3
$ws.Range("D24").Formula = "=SUM((A24/1000)*B24*C24)"
$ws.Range("E24").Value = "Google Play"

# ---------------------------------------------------------------------
# Row 25: grand total of the two price cells above (new formula cell)
# ---------------------------------------------------------------------
$ws.Range("D25").Formula = "=SUM(D23:D24)"
# E25 ("I alt") is unchanged.

# ---------------------------------------------------------------------
# Row 33: new hard-coded total (was =SUM(B29:L29))
# ---------------------------------------------------------------------
$ws.Range("A33").Value = 2000000
# B33, C33, D33's formula and E33 are unchanged.

# ---------------------------------------------------------------------
# Row 34: new hard-coded total (was =SUM(B30:L30))
# ---------------------------------------------------------------------
$ws.Range("A34").Value = 2000000
# B34, C34, D34's formula and E34 are unchanged.

# Row 35's D35 (=SUM(D33:D34)) is an existing formula and recalculates
# automatically.

# ---------------------------------------------------------------------
# sheetView: scroll position + selection
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("A32:E35").Select()

$excel.Calculate()
